$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 54-66 ---
$ws = $wb.Worksheets.Item("PIR")
$dateSrc = $ws.Range("A2")
$dateSrc.Copy($ws.Range("A54"))
$ws.Range("B54").Value = '16:14:51'
$ws.Range("C54").Value = '16:00'
$ws.Range("D54").Value = 'Bathroom'
$ws.Range("E54").Value = 'No Motion'
$ws.Range("F54").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A55"))
$ws.Range("B55").Value = '16:14:54'
$ws.Range("C55").Value = '16:00'
$ws.Range("D55").Value = 'Bathroom'
$ws.Range("E55").Value = 'No Motion'
$ws.Range("F55").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A56"))
$ws.Range("B56").Value = '16:14:59'
$ws.Range("C56").Value = '16:00'
$ws.Range("D56").Value = 'Bathroom'
$ws.Range("E56").Value = 'No Motion'
$ws.Range("F56").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A57"))
$ws.Range("B57").Value = '16:15:04'
$ws.Range("C57").Value = '16:00'
$ws.Range("D57").Value = 'Bathroom'
$ws.Range("E57").Value = 'No Motion'
$ws.Range("F57").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A58"))
$ws.Range("B58").Value = '16:15:09'
$ws.Range("C58").Value = '16:00'
$ws.Range("D58").Value = 'Bathroom'
$ws.Range("E58").Value = 'No Motion'
$ws.Range("F58").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A59"))
$ws.Range("B59").Value = '16:15:14'
$ws.Range("C59").Value = '16:00'
$ws.Range("D59").Value = 'Bathroom'
$ws.Range("E59").Value = 'No Motion'
$ws.Range("F59").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A60"))
$ws.Range("B60").Value = '16:15:19'
$ws.Range("C60").Value = '16:00'
$ws.Range("D60").Value = 'Bathroom'
$ws.Range("E60").Value = 'No Motion'
$ws.Range("F60").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A61"))
$ws.Range("B61").Value = '16:15:24'
$ws.Range("C61").Value = '16:00'
$ws.Range("D61").Value = 'Bathroom'
$ws.Range("E61").Value = 'No Motion'
$ws.Range("F61").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A62"))
$ws.Range("B62").Value = '16:15:29'
$ws.Range("C62").Value = '16:00'
$ws.Range("D62").Value = 'Bathroom'
$ws.Range("E62").Value = 'No Motion'
$ws.Range("F62").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A63"))
$ws.Range("B63").Value = '16:15:34'
$ws.Range("C63").Value = '16:00'
$ws.Range("D63").Value = 'Bathroom'
$ws.Range("E63").Value = 'No Motion'
$ws.Range("F63").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A64"))
$ws.Range("B64").Value = '16:15:39'
$ws.Range("C64").Value = '16:00'
$ws.Range("D64").Value = 'Bathroom'
$ws.Range("E64").Value = 'No Motion'
$ws.Range("F64").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A65"))
$ws.Range("B65").Value = '16:15:44'
$ws.Range("C65").Value = '16:00'
$ws.Range("D65").Value = 'Bathroom'
$ws.Range("E65").Value = 'No Motion'
$ws.Range("F65").Value = 'Inactive'
$dateSrc.Copy($ws.Range("A66"))
$ws.Range("B66").Value = '16:15:49'
$ws.Range("C66").Value = '16:00'
$ws.Range("D66").Value = 'Bathroom'
$ws.Range("E66").Value = 'No Motion'
$ws.Range("F66").Value = 'Inactive'

# --- Humidity sheet: append rows 52-65 ---
$ws = $wb.Worksheets.Item("Humidity")
$dateSrc = $ws.Range("A2")
$dateSrc.Copy($ws.Range("A52"))
$ws.Range("B52").Value = '16:14:51'
$ws.Range("C52").Value = '16:00'
$ws.Range("D52").Value = 'Bathroom'
$ws.Range("E32").Copy($ws.Range("E52"))
$ws.Range("F52").Value = 'Active'
$dateSrc.Copy($ws.Range("A53"))
$ws.Range("B53").Value = '16:14:54'
$ws.Range("C53").Value = '16:00'
$ws.Range("D53").Value = 'Bathroom'
$ws.Range("E35").Copy($ws.Range("E53"))
$ws.Range("F53").Value = 'Active'
$dateSrc.Copy($ws.Range("A54"))
$ws.Range("B54").Value = '16:14:58'
$ws.Range("C54").Value = '16:00'
$ws.Range("D54").Value = 'Bathroom'
$ws.Range("E18").Copy($ws.Range("E54"))
$ws.Range("F54").Value = 'Active'
$dateSrc.Copy($ws.Range("A55"))
$ws.Range("B55").Value = '16:15:02'
$ws.Range("C55").Value = '16:00'
$ws.Range("D55").Value = 'Bathroom'
$ws.Range("E32").Copy($ws.Range("E55"))
$ws.Range("F55").Value = 'Active'
$dateSrc.Copy($ws.Range("A56"))
$ws.Range("B56").Value = '16:15:06'
$ws.Range("C56").Value = '16:00'
$ws.Range("D56").Value = 'Bathroom'
$ws.Range("E18").Copy($ws.Range("E56"))
$ws.Range("F56").Value = 'Active'
$dateSrc.Copy($ws.Range("A57"))
$ws.Range("B57").Value = '16:15:10'
$ws.Range("C57").Value = '16:00'
$ws.Range("D57").Value = 'Bathroom'
$ws.Range("E32").Copy($ws.Range("E57"))
$ws.Range("F57").Value = 'Active'
$dateSrc.Copy($ws.Range("A58"))
$ws.Range("B58").Value = '16:15:18'
$ws.Range("C58").Value = '16:00'
$ws.Range("D58").Value = 'Bathroom'
$ws.Range("E35").Copy($ws.Range("E58"))
$ws.Range("F58").Value = 'Active'
$dateSrc.Copy($ws.Range("A59"))
$ws.Range("B59").Value = '16:15:22'
$ws.Range("C59").Value = '16:00'
$ws.Range("D59").Value = 'Bathroom'
$ws.Range("E16").Copy($ws.Range("E59"))
$ws.Range("F59").Value = 'Active'
$dateSrc.Copy($ws.Range("A60"))
$ws.Range("B60").Value = '16:15:26'
$ws.Range("C60").Value = '16:00'
$ws.Range("D60").Value = 'Bathroom'
$ws.Range("E18").Copy($ws.Range("E60"))
$ws.Range("F60").Value = 'Active'
$dateSrc.Copy($ws.Range("A61"))
$ws.Range("B61").Value = '16:15:30'
$ws.Range("C61").Value = '16:00'
$ws.Range("D61").Value = 'Bathroom'
$ws.Range("E32").Copy($ws.Range("E61"))
$ws.Range("F61").Value = 'Active'
$dateSrc.Copy($ws.Range("A62"))
$ws.Range("B62").Value = '16:15:38'
$ws.Range("C62").Value = '16:00'
$ws.Range("D62").Value = 'Bathroom'
$ws.Range("E18").Copy($ws.Range("E62"))
$ws.Range("F62").Value = 'Active'
$dateSrc.Copy($ws.Range("A63"))
$ws.Range("B63").Value = '16:15:42'
$ws.Range("C63").Value = '16:00'
$ws.Range("D63").Value = 'Bathroom'
$ws.Range("E16").Copy($ws.Range("E63"))
$ws.Range("F63").Value = 'Active'
$dateSrc.Copy($ws.Range("A64"))
$ws.Range("B64").Value = '16:15:46'
$ws.Range("C64").Value = '16:00'
$ws.Range("D64").Value = 'Bathroom'
$ws.Range("E18").Copy($ws.Range("E64"))
$ws.Range("F64").Value = 'Active'
$dateSrc.Copy($ws.Range("A65"))
$ws.Range("B65").Value = '16:15:50'
$ws.Range("C65").Value = '16:00'
$ws.Range("D65").Value = 'Bathroom'
$ws.Range("E16").Copy($ws.Range("E65"))
$ws.Range("F65").Value = 'Active'

# --- Temperature sheet: append rows 52-65 ---
$ws = $wb.Worksheets.Item("Temperature")
$dateSrc = $ws.Range("A2")
$dateSrc.Copy($ws.Range("A52"))
$ws.Range("B52").Value = '16:14:51'
$ws.Range("C52").Value = '16:00'
$ws.Range("D52").Value = 'Bathroom'
$ws.Range("E52").Value = '22.8C'
$ws.Range("F52").Value = 'Active'
$dateSrc.Copy($ws.Range("A53"))
$ws.Range("B53").Value = '16:14:54'
$ws.Range("C53").Value = '16:00'
$ws.Range("D53").Value = 'Bathroom'
$ws.Range("E53").Value = '22.8C'
$ws.Range("F53").Value = 'Active'
$dateSrc.Copy($ws.Range("A54"))
$ws.Range("B54").Value = '16:14:58'
$ws.Range("C54").Value = '16:00'
$ws.Range("D54").Value = 'Bathroom'
$ws.Range("E54").Value = '22.8C'
$ws.Range("F54").Value = 'Active'
$dateSrc.Copy($ws.Range("A55"))
$ws.Range("B55").Value = '16:15:02'
$ws.Range("C55").Value = '16:00'
$ws.Range("D55").Value = 'Bathroom'
$ws.Range("E55").Value = '22.8C'
$ws.Range("F55").Value = 'Active'
$dateSrc.Copy($ws.Range("A56"))
$ws.Range("B56").Value = '16:15:06'
$ws.Range("C56").Value = '16:00'
$ws.Range("D56").Value = 'Bathroom'
$ws.Range("E56").Value = '22.8C'
$ws.Range("F56").Value = 'Active'
$dateSrc.Copy($ws.Range("A57"))
$ws.Range("B57").Value = '16:15:10'
$ws.Range("C57").Value = '16:00'
$ws.Range("D57").Value = 'Bathroom'
$ws.Range("E57").Value = '22.8C'
$ws.Range("F57").Value = 'Active'
$dateSrc.Copy($ws.Range("A58"))
$ws.Range("B58").Value = '16:15:19'
$ws.Range("C58").Value = '16:00'
$ws.Range("D58").Value = 'Bathroom'
$ws.Range("E58").Value = '22.7C'
$ws.Range("F58").Value = 'Active'
$dateSrc.Copy($ws.Range("A59"))
$ws.Range("B59").Value = '16:15:22'
$ws.Range("C59").Value = '16:00'
$ws.Range("D59").Value = 'Bathroom'
$ws.Range("E59").Value = '22.8C'
$ws.Range("F59").Value = 'Active'
$dateSrc.Copy($ws.Range("A60"))
$ws.Range("B60").Value = '16:15:27'
$ws.Range("C60").Value = '16:00'
$ws.Range("D60").Value = 'Bathroom'
$ws.Range("E60").Value = '22.8C'
$ws.Range("F60").Value = 'Active'
$dateSrc.Copy($ws.Range("A61"))
$ws.Range("B61").Value = '16:15:31'
$ws.Range("C61").Value = '16:00'
$ws.Range("D61").Value = 'Bathroom'
$ws.Range("E61").Value = '22.7C'
$ws.Range("F61").Value = 'Active'
$dateSrc.Copy($ws.Range("A62"))
$ws.Range("B62").Value = '16:15:39'
$ws.Range("C62").Value = '16:00'
$ws.Range("D62").Value = 'Bathroom'
$ws.Range("E62").Value = '22.8C'
$ws.Range("F62").Value = 'Active'
$dateSrc.Copy($ws.Range("A63"))
$ws.Range("B63").Value = '16:15:43'
$ws.Range("C63").Value = '16:00'
$ws.Range("D63").Value = 'Bathroom'
$ws.Range("E63").Value = '22.8C'
$ws.Range("F63").Value = 'Active'
$dateSrc.Copy($ws.Range("A64"))
$ws.Range("B64").Value = '16:15:47'
$ws.Range("C64").Value = '16:00'
$ws.Range("D64").Value = 'Bathroom'
$ws.Range("E64").Value = '22.8C'
$ws.Range("F64").Value = 'Active'
$dateSrc.Copy($ws.Range("A65"))
$ws.Range("B65").Value = '16:15:51'
$ws.Range("C65").Value = '16:00'
$ws.Range("D65").Value = 'Bathroom'
$ws.Range("E65").Value = '22.8C'
$ws.Range("F65").Value = 'Active'
